# The workbook has 4 language sheets (MR, HI, TA, EN) that each repeat the
# same 4-column header row (Sr. No. / Code / Name (Roman) / Name (Orig)).
# The author renamed the shared "Name (Orig)" header to "Name (Original)"
# (Excel re-uses/reorders the shared-string table for this across every
# sheet that referenced it), widened the HI sheet's 4th column so the new,
# longer header still fits, and finished up with the EN tab selected/active
# with the cursor sitting on D2.

$wb = $excel.ActiveWorkbook

$sheetNames = @("MR", "HI", "TA", "EN")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("D1").Value = "Name (Original)"
    $ws.Range("D2").Select()
}

# HI's column D needs to widen from 10.90625 to 14 to fit the longer header.
$wb.Worksheets.Item("HI").Columns.Item(4).ColumnWidth = 13.166666666666666

# EN becomes the active sheet/tab (was TA before).
$wb.Worksheets.Item("EN").Activate()
